$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.001144417328760028
$ws.Range("C3").Value = 0.005653401836752892
$ws.Range("C4").Value = 0.008702149614691734
$ws.Range("C5").Value = 0.01114099565893412
$ws.Range("C6").Value = 0.01293173339217901
$ws.Range("C7").Value = 0.01401621662080288
$ws.Range("C8").Value = 0.0144551582634449
$ws.Range("C9").Value = 0.01450350880622864
$ws.Range("C10").Value = 0.01437765080481768
$ws.Range("C11").Value = 0.01428111549466848
